$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (row 12) down into the new row 13
$ws.Range("A12:N12").Copy()
$ws.Range("A13:N13").PasteSpecial(-4122)

$row = 13
$ws.Cells.Item($row, 1).Value = 42620.891192129631
$ws.Cells.Item($row, 2).Value = -4
$ws.Cells.Item($row, 3).Value = 52
$ws.Cells.Item($row, 4).Value = 43
$ws.Cells.Item($row, 5).Value = 52
$ws.Cells.Item($row, 6).Value = 50
$ws.Cells.Item($row, 7).Value = 34676
$ws.Cells.Item($row, 8).Value = 22975
$ws.Cells.Item($row, 9).Value = 1179
$ws.Cells.Item($row, 10).Value = 261
$ws.Cells.Item($row, 11).Value = 216
$ws.Cells.Item($row, 12).Value = 1
$ws.Cells.Item($row, 13).Value = 1
$ws.Cells.Item($row, 14).Value = "Named"
